$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "imię"
$ws.Range("D2").Value = "Zuzanna błaszczak"
$ws.Range("D3").Value = "ktokolwiek ktokolwiek"

$ws.Range("D8").Select()
